$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume/1h%) to match the new source
# snapshot. Price cells are stored as text in this sheet (e.g. "88.70",
# "0.713"), so numeric-looking replacements are entered with a leading
# apostrophe to force text entry and avoid Excel coercing them into numbers
# (which would drop significant trailing zeros). Rows 36/37 (OKB <->
# Bittensor) and 50/51 (ApeXProtocol <-> Maker) swap rank order, so all of
# name/link/price/volume are rewritten for those four rows.

$ws.Range("D2").Value = '68.129.21'
$ws.Range("E2").Value = '  -3.83%  '
$ws.Range("D3").Value = '3.687.72'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'594.65"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'180.69"
$ws.Range("E6").Value = '  +8.81%  '
$ws.Range("D7").Value = '3.678.63'
$ws.Range("E7").Value = '  -4.66%  '
$ws.Range("E8").Value = '  -6.84%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = "'0.713"
$ws.Range("E10").Value = '  -5.16%  '
$ws.Range("E11").Value = '  -7.69%  '
$ws.Range("D12").Value = "'56.14"
$ws.Range("E12").Value = '  +5.49%  '
$ws.Range("D13").Value = "'0.0000291"
$ws.Range("E13").Value = '  -9.44%  '
$ws.Range("E14").Value = '  -9.26%  '
$ws.Range("D15").Value = '4.270.50'
$ws.Range("E15").Value = '  -4.86%  '
$ws.Range("D16").Value = '3.682.33'
$ws.Range("E16").Value = '  -4.65%  '
$ws.Range("D17").Value = "'19.26"
$ws.Range("E17").Value = '  -9.80%  '
$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = '  -7.76%  '
$ws.Range("E20").Value = '  -7.62%  '
$ws.Range("D21").Value = '68.000.39'
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").Value = "'408.13"
$ws.Range("E22").Value = '  -6.53%  '
$ws.Range("E23").Value = '  -3.23%  '
$ws.Range("D24").Value = "'88.70"
$ws.Range("E24").Value = '  -5.80%  '
$ws.Range("E25").Value = '  -7.78%  '
$ws.Range("D26").Value = "'12.77"
$ws.Range("E26").Value = '  -7.85%  '
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = '  -3.30%  '
$ws.Range("D28").Value = "'3.86"
$ws.Range("E28").Value = '  -5.01%  '
$ws.Range("D29").Value = "'6.05"
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("D30").Value = "'9.40"
$ws.Range("E30").Value = '  -9.52%  '
$ws.Range("E31").Value = '  -6.81%  '
$ws.Range("D32").Value = "'7.23"
$ws.Range("E32").Value = '  -9.57%  '
$ws.Range("D33").Value = "'12.42"
$ws.Range("E33").Value = '  -8.25%  '
$ws.Range("E34").Value = '  -6.69%  '
$ws.Range("D35").Value = "'43.36"
$ws.Range("E35").Value = '  -9.47%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = "'603.18"
$ws.Range("E36").Value = '  -4.91%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = "'64.20"
$ws.Range("E37").Value = '  -8.17%  '
$ws.Range("D38").Value = '0.0₃0884'
$ws.Range("E38").Value = '  -9.96%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E40").Value = '  -6.02%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  -7.31%  '
$ws.Range("D43").Value = "'2.79"
$ws.Range("E43").Value = '  +3.16%  '
$ws.Range("D44").Value = "'2.99"
$ws.Range("E44").Value = '  -9.48%  '
$ws.Range("D45").Value = "'0.0437"
$ws.Range("E45").Value = '  -6.89%  '
$ws.Range("D46").Value = "'2.84"
$ws.Range("E46").Value = '  -12.56%  '
$ws.Range("E47").Value = '  -8.85%  '
$ws.Range("D48").Value = "'2.72"
$ws.Range("E48").Value = '  -4.26%  '
$ws.Range("E49").Value = '  -6.94%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.733.23'
$ws.Range("E50").Value = '  -3.55%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = "'3.14"
$ws.Range("E51").Value = '  -6.05%  '
